$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column J header in row 1 (value 8), copying style from I1
$ws.Cells.Item(1, 9).Copy()
$ws.Cells.Item(1, 10).PasteSpecial(-4122)
$ws.Cells.Item(1, 10).Value = 8

# Row 2
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0.003528581510232887
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0.01909241837299393

# Row 3
$ws.Cells.Item(3, 2).Value = 0.1443066516347239
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0.001411432604093155
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0.2382528127068175
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0.02296624239070285

# Row 4
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 0.1057692307692306
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0.07317073170731707
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0.06159895150720845
$ws.Cells.Item(4, 10).Value = 0

# Row 5
$ws.Cells.Item(5, 2).Value = 0.0541149943630213
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.04446012702893437
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0.02117802779616148
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0.05368013281682354

# Row 6
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0.003528581510232887
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 0.02877697841726622

# Row 7
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0.01736972704714641
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0.1441677588466578
$ws.Cells.Item(7, 10).Value = 0

# Row 8
$ws.Cells.Item(8, 2).Value = 0.08399098083427288
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.01199717713479181
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0.2183984116479159
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 0.05838406198118436

# Row 9
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0.2142266335814727
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0.02328863796753705
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0.01048492791612058
$ws.Cells.Item(9, 10).Value = 0

# Row 10
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0.1038461538461537
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 0.1219512195121951
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 0.1153342070773263
$ws.Cells.Item(10, 10).Value = 0

# Row 11
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.0008271298593879239
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0.01764290755116443
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = 0

# Row 12
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 0.02646815550041358
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 0.0155257586450247
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0.01703800786369594
$ws.Cells.Item(12, 10).Value = 0

# Row 13
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0.009174311926605505
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0.01300498063087992

# Row 14
$ws.Cells.Item(14, 2).Value = 0.0259301014656144
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0.0007057163020465773
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 0.0113447703375761

# Row 15
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0

# Row 16
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 0.01153846153846154
$ws.Cells.Item(16, 4).Value = 0.1819685690653434
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 0.2173606210303458
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0.1284403669724769
$ws.Cells.Item(16, 10).Value = 0.02822357498616495

# Row 17
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0.02767017155506368

# Row 18
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 0.009925558312655087
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0.007863695937090432
$ws.Cells.Item(18, 10).Value = 0

# Row 19
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 0.02399435426958363
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0.01023796347537356

# Row 20
$ws.Cells.Item(20, 2).Value = 0.001691093573844419
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0.006351446718419196
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0.1242390702822359

# Row 21
$ws.Cells.Item(21, 2).Value = 0.02085682074408116
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = 0.001654259718775848
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0.05363443895553988
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0.1098505810736028

# Row 22
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0

# Row 23
$ws.Cells.Item(23, 2).Value = 0.1104847801578357
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0.02117148906139732
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0.170747849106552
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0.1045932484781407

# Row 24
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0.0008271298593879239
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 0.0432503276539974
$ws.Cells.Item(24, 10).Value = 0

# Row 25
$ws.Cells.Item(25, 2).Value = 0
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 0.04549214226633585
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0.04940014114326041
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0

# Row 26
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 0.0008271298593879239
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 0

# Row 27
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 0.01902398676592226
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 0.01129146083274524
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0.01965923984272608
$ws.Cells.Item(27, 10).Value = 0

# Row 28
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0.01270289343683839
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0.02628666297731049

# Row 29
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 0.03881439661256175
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0.002490315439955728

# Row 30
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 0.02233250620347396
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0.05645730416372619
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0.001310615989515072
$ws.Cells.Item(30, 10).Value = 0

# Row 31
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 0.07526881720430104
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0.02883355176933158
$ws.Cells.Item(31, 10).Value = 0

# Row 32
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 0.6365384615384595
$ws.Cells.Item(32, 4).Value = 0.06534325889164602
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 0.8048780487804883
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0.1572739187418084
$ws.Cells.Item(32, 10).Value = 0

# Row 33
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 0.142307692307692
$ws.Cells.Item(33, 4).Value = 0.1248966087675761
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 0.004234297812279464
$ws.Cells.Item(33, 7).Value = 0
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0.2005242463958057
$ws.Cells.Item(33, 10).Value = 0

# Row 34
$ws.Cells.Item(34, 2).Value = 0
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0.0007057163020465773
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0.01162147205312674

# Row 35
$ws.Cells.Item(35, 2).Value = 0
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 0.04383788254756
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0.002117148906139732
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0.007863695937090432
$ws.Cells.Item(35, 10).Value = 0

# Delete rows 36-40 (Joint regime area rows removed)
$ws.Range("A36:J40").EntireRow.Delete()